$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.001" or
# "10.60" are not auto-converted/rounded to numbers by the smart typing in Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.969.27"
$ws.Range("E2").Value = "  +5.39%  "
$ws.Range("D3").Value = "1.914.23"
$ws.Range("E3").Value = "  +4.67%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "338.91"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "0.4731"
$ws.Range("E7").Value = "  +3.19%  "
$ws.Range("D8").Value = "0.4047"
$ws.Range("E8").Value = "  +6.23%  "
$ws.Range("D9").Value = "48.12"
$ws.Range("D10").Value = "0.08171"
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("D11").Value = "1.030"
$ws.Range("E11").Value = "  +6.05%  "
$ws.Range("D12").Value = "22.46"
$ws.Range("E12").Value = "  +6.59%  "
$ws.Range("D13").Value = "1.899.75"
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").Value = "6.097"
$ws.Range("E14").Value = "  +3.56%  "
$ws.Range("D15").Value = "7.376"
$ws.Range("E15").Value = "  +4.40%  "
$ws.Range("D16").Value = "91.52"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "0.00001054"
$ws.Range("E18").Value = "  +2.75%  "
$ws.Range("D19").Value = "0.06647"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "17.87"
$ws.Range("E20").Value = "  +4.38%  "
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "28.989.17"
$ws.Range("E22").Value = "  +5.55%  "
$ws.Range("D23").Value = "5.563"
$ws.Range("E23").Value = "  +4.22%  "
$ws.Range("D24").Value = "11.20"
$ws.Range("E24").Value = "  +3.52%  "
$ws.Range("D25").Value = "2.271"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").Value = "2.129.67"
$ws.Range("E26").Value = "  +3.89%  "
$ws.Range("D27").Value = "160.86"
$ws.Range("E27").Value = "  +3.45%  "
$ws.Range("D28").Value = "20.07"
$ws.Range("E28").Value = "  +3.51%  "
$ws.Range("D29").Value = "2.181"
$ws.Range("E29").Value = "  +5.82%  "
$ws.Range("D30").Value = "5.526"
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D31").Value = "120.98"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "1.015"
$ws.Range("E32").Value = "  +7.66%  "
$ws.Range("D33").Value = "0.09587"
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.661"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.430"
$ws.Range("E35").Value = "  +7.62%  "
$ws.Range("D36").Value = "5.426"
$ws.Range("E36").Value = "  +3.27%  "
$ws.Range("D37").Value = "0.06227"
$ws.Range("E37").Value = "  +4.90%  "
$ws.Range("D38").Value = "0.02289"
$ws.Range("E38").Value = "  +4.64%  "
$ws.Range("D39").Value = "8.669"
$ws.Range("E39").Value = "  +7.53%  "
$ws.Range("D40").Value = "1.202"
$ws.Range("E40").Value = "  +4.99%  "
$ws.Range("D41").Value = "0.6043"
$ws.Range("E41").Value = "  +4.63%  "
$ws.Range("D42").Value = "10.60"
$ws.Range("E42").Value = "  +6.06%  "
$ws.Range("D43").Value = "0.1904"
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("D44").Value = "0.9986"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "1.272"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "0.5643"
$ws.Range("E46").Value = "  +3.48%  "
$ws.Range("D47").Value = "12.37"
$ws.Range("E47").Value = "  +3.53%  "
$ws.Range("D48").Value = "1.982"
$ws.Range("E48").Value = "  +5.99%  "
$ws.Range("D49").Value = "0.07314"
$ws.Range("E49").Value = "  +10.71%  "
$ws.Range("D50").Value = "2.151"
$ws.Range("E50").Value = "  +18.86%  "
$ws.Range("D51").Value = "113.29"
$ws.Range("E51").Value = "  +2.03%  "

# Restore the default (unstyled) cell style on column D so formatting
# matches the original workbook (no explicit number format retained).
$ws.Range("D2:D51").Style = "Normal"

